$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jawaban")

# Remove the old tail content (rows 5-8 and columns E-H) so the used
# range shrinks back down to A1:D4, mirroring the diff.
$ws.Range("A5:H8").ClearContents()
$ws.Range("E1:H4").ClearContents()

# Overwrite the staircase of single-letter cells with the new spelling
# ("kodepython") in place of the old ("Purwadhika StartupandCodingSchool@BSD").
$ws.Range("A1").Value = "k"

$ws.Range("A2").Value = "o"
$ws.Range("B2").Value = "d"

$ws.Range("A3").Value = "e"
$ws.Range("B3").Value = "p"
$ws.Range("C3").Value = "y"

$ws.Range("A4").Value = "t"
$ws.Range("B4").Value = "h"
$ws.Range("C4").Value = "o"
$ws.Range("D4").Value = "n"
